# feat: add 2022-Q4 data
#
# 1. Duplicate the existing "2022-Q3" detail sheet (to inherit its exact
#    layout/formatting), rename the duplicate to "2022-Q4", position it
#    right before "2022-Q3", and overwrite its contents with the new
#    2022-Q4 fund-holding data (7 rows instead of the old 9).
# 2. Insert a new leading row into the "总计" (summary) sheet for the
#    2022-Q4 totals, shifting the existing quarters down by one and
#    renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q4" detail sheet
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($wsQ3)
$wsQ4 = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ4.Name = "2022-Q4"

$q4Data = @(
    @("013958", "华商鑫选回报一年持有混合A", "5.33", "90.81", "1.92", "0.1023", 5),
    @("001170", "泰达宏利复兴伟业灵活配置混合", "1.64", "91.15", "4.78", "0.0784", 7),
    @("013959", "华商鑫选回报一年持有混合C", "1.22", "90.81", "1.92", "0.0234", 5),
    @("012216", "红塔红土盛利混合A", "2.23", "50.74", "1.01", "0.0225", 9),
    @("013733", "红塔红土盛丰混合A", "0.40", "61.27", "4.50", "0.0180", 2),
    @("013734", "红塔红土盛丰混合C", "0.09", "61.27", "4.50", "0.0040", 2),
    @("012217", "红塔红土盛利混合C", "0.08", "50.74", "1.01", "0.0008", 9)
)

# The copied "2022-Q3" sheet has 9 data rows (rows 2-10); the new
# "2022-Q4" sheet only needs 7 (rows 2-8), so drop the two extra rows.
$wsQ4.Rows.Item(9).Delete()
$wsQ4.Rows.Item(9).Delete()

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = $i + 2
    $rec = $q4Data[$i]

    $wsQ4.Cells.Item($row, 1).Value = $i

    # Columns B-G (code, name, size, position, ratio, held value) are
    # plain text in this workbook, even though several look numeric
    # (e.g. "013958", "0.40") - force text storage without leaving any
    # explicit cell style behind, matching the original sheets' cells.
    for ($c = 2; $c -le 7; $c++) {
        $cell = $wsQ4.Cells.Item($row, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$c - 2]
        $cell.ClearFormats()
    }

    # Column H (rank) is numeric.
    $wsQ4.Cells.Item($row, 8).Value = $rec[6]
}

# ---------------------------------------------------------------------
# Step 2: add the 2022-Q4 row to the "总计" summary sheet
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# Shift existing rows 2-5 down to 3-6 (bottom-up copy to avoid clobbering).
for ($r = 5; $r -ge 2; $r--) {
    $dest = $r + 1
    $wsTotal.Cells.Item($r, 1).Copy($wsTotal.Cells.Item($dest, 1))
    $wsTotal.Cells.Item($r, 2).Copy($wsTotal.Cells.Item($dest, 2))
    $wsTotal.Cells.Item($r, 3).Copy($wsTotal.Cells.Item($dest, 3))
    $wsTotal.Cells.Item($r, 4).Copy($wsTotal.Cells.Item($dest, 4))
}

# New first data row: 2022-Q4
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q4"
$wsTotal.Cells.Item(2, 3).Value = 7
$wsTotal.Cells.Item(2, 4).Value = 0.25

# Column A is a plain 0-based row index, not a copied value - renumber it
# sequentially across the now-6-row table (rows 2-6 => 0,1,2,3,4).
for ($r = 2; $r -le 6; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}
